$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (I1) and IF (J1), copying H1's format (bold/border/centered)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-18: column I = 1 (constant), column J = copy of column H
for ($r = 2; $r -le 18; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
